$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("H62").Value = 7054.1665
$ws.Range("I62").Value = 6276.6665
$ws.Range("K62").Value = 6276.6665
$ws.Range("M62").Value = -5652.6665
$ws.Range("H65").Value = 7054.1665
$ws.Range("I65").Value = 6276.6665
$ws.Range("K65").Value = 31383.3325
$ws.Range("M65").Value = -28263.3325
$ws.Range("H70").Value = 3107.3333
$ws.Range("I70").Value = 2689
$ws.Range("J70").Value = 3316.5
$ws.Range("K70").Value = 8067
$ws.Range("L70").Value = 9949.5
$ws.Range("M70").Value = -7797
$ws.Range("N70").Value = -10489.5
$ws.Range("H73").Value = 3107.3333
$ws.Range("I73").Value = 2689
$ws.Range("J73").Value = 3316.5
$ws.Range("K73").Value = 8067
$ws.Range("L73").Value = 9949.5
$ws.Range("M73").Value = -7131
$ws.Range("N73").Value = -11821.5
$ws.Range("H80").Value = 2182.1072
$ws.Range("I80").Value = 2638.5
$ws.Range("J80").Value = 1360.6
$ws.Range("K80").Value = 7915.5
$ws.Range("L80").Value = 4081.8
$ws.Range("M80").Value = -6917.5
$ws.Range("N80").Value = -6077.799999999999
$ws.Range("H83").Value = 2182.1072
$ws.Range("I83").Value = 2638.5
$ws.Range("J83").Value = 1360.6
$ws.Range("K83").Value = 23746.5
$ws.Range("L83").Value = 12245.4
$ws.Range("M83").Value = -18754.5
$ws.Range("N83").Value = -22229.4
$ws.Range("H129").Value = 15108.5625
$ws.Range("J129").Value = 29063.875
$ws.Range("L129").Value = 87191.625
$ws.Range("N129").Value = -97191.625

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1077.9615
$ws.Range("I2").Value = 779.55554
$ws.Range("K2").Value = 779.55554
$ws.Range("M2").Value = -666.55554
$ws.Range("H32").Value = 12508687
$ws.Range("I32").Value = 17244860
$ws.Range("K32").Value = 17244860
$ws.Range("M32").Value = -17244573
$ws.Range("H61").Value = 37041176
$ws.Range("I61").Value = 66669588
$ws.Range("K61").Value = 66669588
$ws.Range("M61").Value = -66669376
$ws.Range("H88").Value = 2787
$ws.Range("J88").Value = 2871.75
$ws.Range("L88").Value = 2871.75
$ws.Range("N88").Value = -3683.75
$ws.Range("H91").Value = 2787
$ws.Range("J91").Value = 2871.75
$ws.Range("L91").Value = 2871.75
$ws.Range("N91").Value = -5679.75
$ws.Range("H97").Value = 1987.7693
$ws.Range("I97").Value = 568.25
$ws.Range("J97").Value = 4259
$ws.Range("K97").Value = 568.25
$ws.Range("L97").Value = 4259
$ws.Range("M97").Value = -72.25
$ws.Range("N97").Value = -5251
$ws.Range("H116").Value = 1077.9615
$ws.Range("I116").Value = 779.55554
$ws.Range("K116").Value = 779.55554
$ws.Range("M116").Value = 1514.44446
$ws.Range("H122").Value = 2284.2058
$ws.Range("I122").Value = 1269.7727
$ws.Range("J122").Value = 4144
$ws.Range("K122").Value = 3809.3181
$ws.Range("L122").Value = 12432
$ws.Range("M122").Value = -1359.3181
$ws.Range("N122").Value = -17332
$ws.Range("H124").Value = 45140.715
$ws.Range("J124").Value = 45140.715
$ws.Range("L124").Value = 45140.715
$ws.Range("N124").Value = -54960.715
$ws.Range("H132").Value = 45466716
$ws.Range("J132").Value = 250002320
$ws.Range("L132").Value = 750006960
$ws.Range("N132").Value = -750012020
$ws.Range("H136").Value = 37041176
$ws.Range("I136").Value = 66669588
$ws.Range("K136").Value = 200008764
$ws.Range("M136").Value = -200006214

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1077.9615
$ws.Range("I3").Value = 779.55554
$ws.Range("K3").Value = 779.55554
$ws.Range("M3").Value = -665.55554
$ws.Range("H75").Value = 20082.1
$ws.Range("I75").Value = 5403
$ws.Range("J75").Value = 54333.332
$ws.Range("K75").Value = 5403
$ws.Range("L75").Value = 54333.332
$ws.Range("M75").Value = -4467
$ws.Range("N75").Value = -56205.332
$ws.Range("H78").Value = 20082.1
$ws.Range("I78").Value = 5403
$ws.Range("J78").Value = 54333.332
$ws.Range("K78").Value = 16209
$ws.Range("L78").Value = 162999.996
$ws.Range("M78").Value = -11529
$ws.Range("N78").Value = -172359.996
$ws.Range("H86").Value = 12789.685
$ws.Range("I86").Value = 2158.3333
$ws.Range("J86").Value = 52657.25
$ws.Range("K86").Value = 2158.3333
$ws.Range("L86").Value = 52657.25
$ws.Range("M86").Value = -1035.3333
$ws.Range("N86").Value = -54903.25
$ws.Range("H89").Value = 12789.685
$ws.Range("I89").Value = 2158.3333
$ws.Range("J89").Value = 52657.25
$ws.Range("K89").Value = 10791.6665
$ws.Range("L89").Value = 263286.25
$ws.Range("M89").Value = -5175.666499999999
$ws.Range("N89").Value = -274518.25
$ws.Range("H92").Value = 134720.75
$ws.Range("J92").Value = 134720.75
$ws.Range("L92").Value = 134720.75
$ws.Range("N92").Value = -139712.75
$ws.Range("H134").Value = 5261.92
$ws.Range("I134").Value = 4973.9414
$ws.Range("J134").Value = 5873.875
$ws.Range("K134").Value = 14921.8242
$ws.Range("L134").Value = 17621.625
$ws.Range("M134").Value = -12386.8242
$ws.Range("N134").Value = -22691.625

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5253.476
$ws.Range("I22").Value = 7050.7334
$ws.Range("J22").Value = 760.3333
$ws.Range("K22").Value = 7050.7334
$ws.Range("L22").Value = 760.3333
$ws.Range("M22").Value = -6700.7334
$ws.Range("N22").Value = -1460.3333
$ws.Range("H31").Value = 21282704
$ws.Range("I31").Value = 5788.4873
$ws.Range("K31").Value = 5788.4873
$ws.Range("M31").Value = -5493.4873
$ws.Range("H34").Value = 21282704
$ws.Range("I34").Value = 5788.4873
$ws.Range("K34").Value = 5788.4873
$ws.Range("M34").Value = -5586.4873
$ws.Range("H86").Value = 3850.4
$ws.Range("I86").Value = 3100.2
$ws.Range("J86").Value = 4600.6
$ws.Range("K86").Value = 3100.2
$ws.Range("L86").Value = 4600.6
$ws.Range("M86").Value = -1977.2
$ws.Range("N86").Value = -6846.6
$ws.Range("H89").Value = 3850.4
$ws.Range("I89").Value = 3100.2
$ws.Range("J89").Value = 4600.6
$ws.Range("K89").Value = 15501
$ws.Range("L89").Value = 23003
$ws.Range("M89").Value = -9885
$ws.Range("N89").Value = -34235
$ws.Range("H99").Value = 4262.9375
$ws.Range("I99").Value = 4237.364
$ws.Range("K99").Value = 4237.364
$ws.Range("M99").Value = -2739.364
$ws.Range("H105").Value = 8096.5293
$ws.Range("I105").Value = 2032.1111
$ws.Range("K105").Value = 2032.1111
$ws.Range("M105").Value = -285.1111000000001
$ws.Range("H122").Value = 2197.3
$ws.Range("I122").Value = 2293.3333
$ws.Range("K122").Value = 6879.999899999999
$ws.Range("M122").Value = -4429.999899999999
$ws.Range("H126").Value = 4262.9375
$ws.Range("I126").Value = 4237.364
$ws.Range("K126").Value = 12712.092
$ws.Range("M126").Value = -10242.092

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 225777
$ws.Range("I70").Value = 403199.6
$ws.Range("K70").Value = 403199.6
$ws.Range("M70").Value = -402929.6
$ws.Range("H73").Value = 225777
$ws.Range("I73").Value = 403199.6
$ws.Range("K73").Value = 403199.6
$ws.Range("M73").Value = -402263.6
$ws.Range("H98").Value = 202818
$ws.Range("J98").Value = 202818
$ws.Range("L98").Value = 202818
$ws.Range("N98").Value = -208808
$ws.Range("H122").Value = 2934
$ws.Range("I122").Value = 2864.4
$ws.Range("K122").Value = 8593.200000000001
$ws.Range("M122").Value = -6143.200000000001

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1562.8125
$ws.Range("I93").Value = 1311
$ws.Range("J93").Value = 1982.5
$ws.Range("K93").Value = 1311
$ws.Range("L93").Value = 1982.5
$ws.Range("M93").Value = -63
$ws.Range("N93").Value = -4478.5
$ws.Range("H122").Value = 4292.6577
$ws.Range("I122").Value = 3874.8518
$ws.Range("J122").Value = 5318.1816
$ws.Range("K122").Value = 11624.5554
$ws.Range("L122").Value = 15954.5448
$ws.Range("M122").Value = -9174.555399999999
$ws.Range("N122").Value = -20854.5448
$ws.Range("H132").Value = 52633324
$ws.Range("I132").Value = 1686.24
$ws.Range("J132").Value = 153848020
$ws.Range("K132").Value = 5058.72
$ws.Range("L132").Value = 461544060
$ws.Range("M132").Value = -2528.72
$ws.Range("N132").Value = -461549120

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7892.25
$ws.Range("I96").Value = 7231.3335
$ws.Range("J96").Value = 8288.799999999999
$ws.Range("K96").Value = 7231.3335
$ws.Range("L96").Value = 8288.799999999999
$ws.Range("M96").Value = -5858.3335
$ws.Range("N96").Value = -11034.8
$ws.Range("H103").Value = 42416.832
$ws.Range("J103").Value = 42416.832
$ws.Range("L103").Value = 42416.832
$ws.Range("N103").Value = -44760.832
$ws.Range("H109").Value = 75000
$ws.Range("J109").Value = 74999
$ws.Range("L109").Value = 74999
$ws.Range("H122").Value = 37038990
$ws.Range("I122").Value = 52633468
$ws.Range("K122").Value = 157900404
$ws.Range("M122").Value = -157897954
$ws.Range("H123").Value = 78000
$ws.Range("J123").Value = 78000
$ws.Range("L123").Value = 78000
$ws.Range("N123").Value = -87800
$ws.Range("H132").Value = 4924.51
$ws.Range("I132").Value = 4990.1274
$ws.Range("J132").Value = 4153.5
$ws.Range("K132").Value = 14970.3822
$ws.Range("L132").Value = 12460.5
$ws.Range("M132").Value = -12440.3822
$ws.Range("N132").Value = -17520.5
$ws.Range("H136").Value = 1490.8723
$ws.Range("I136").Value = 1431.0238
$ws.Range("K136").Value = 4293.0714
$ws.Range("M136").Value = -1743.0714

# Cell removals (cells deleted entirely, e.g. due to upstream formula/row change)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N12").ClearContents()

# New cells added
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N109").Value = -77773
